# Refresh the coin price/volume snapshot (cryptos list), matching the
# GitHub Actions commit "Updated cryptos list on Wed May  8 23:45:43 UTC 2024".
# Column D/E values are stored as plain text in the sheet (the source data
# keeps formatted strings such as "61.061.47" or "  -2.29%  " rather than
# numbers), so purely-numeric-looking prices are entered with a leading
# apostrophe - exactly what Excel's UI does to force text entry - to stop
# them from being auto-converted into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.061.47"
$ws.Range("E2").Value = "  -2.29%  "
$ws.Range("D3").Value = "2.971.93"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Formula = "'588.35"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("D6").Formula = "'141.94"
$ws.Range("E6").Value = "  -4.77%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Formula = "'0.517"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "2.970.07"
$ws.Range("E10").Value = "  -5.01%  "
$ws.Range("D11").Formula = "'5.76"
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("D12").Formula = "'0.452"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("D14").Formula = "'33.94"
$ws.Range("E14").Value = "  -4.30%  "
$ws.Range("E15").Value = "  +1.87%  "
$ws.Range("D16").Value = "3.463.63"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "61.056.87"
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").Value = "2.969.67"
$ws.Range("E19").Value = "  -1.56%  "
$ws.Range("D20").Formula = "'447.48"
$ws.Range("E20").Value = "  -5.26%  "
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").Formula = "'0.681"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Formula = "'7.33"
$ws.Range("E23").Value = "  -1.43%  "
$ws.Range("D24").Formula = "'81.22"
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("E26").Value = "  -8.78%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Formula = "'9.89"
$ws.Range("E28").Value = "  -5.39%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  +1.10%  "
$ws.Range("D32").Formula = "'2.05"
$ws.Range("E32").Value = "  -5.92%  "
$ws.Range("D33").Formula = "'27.14"
$ws.Range("E33").Value = "  -0.39%  "
$ws.Range("D34").Formula = "'0.106"
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("B35").Value = "Mantle"
$ws.Range("C35").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D35").Formula = "'1.01"
$ws.Range("E35").Value = "  -3.35%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").Value = "0.0₃0783"
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").Formula = "'5.70"
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("D38").Formula = "'50.19"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Formula = "'2.06"
$ws.Range("E39").Value = "  -4.58%  "
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Formula = "'9.11"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").Formula = "'0.118"
$ws.Range("E41").Value = "  +4.80%  "
$ws.Range("D42").Formula = "'2.75"
$ws.Range("E42").Value = "  -9.04%  "
$ws.Range("D43").Formula = "'389.63"
$ws.Range("E43").Value = "  -7.53%  "
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("D45").Value = "2.691.81"
$ws.Range("E45").Value = "  -4.29%  "
$ws.Range("E46").Value = "  -6.68%  "
$ws.Range("D47").Formula = "'37.17"
$ws.Range("E47").Value = "  -2.63%  "
$ws.Range("D48").Formula = "'131.75"
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("E51").Value = "  -0.34%  "
